# df_resumo_performance_escolar.xlsx - "up - final pandas"
#
# The workbook was re-opened/re-saved in Excel; the vast majority of the
# underlying diff (fileVersion/@rupBuild, the x15ac:absPath, the
# xr:revisionPtr documentId, the workbookView window geometry, the
# theme1.xml name/panose/extended-font-list churn, and the
# x14ac:dyDescent rendering hints) is Excel-build metadata that is not
# reachable from the Excel object model (no COM property exposes any of
# it - it is stamped by the Excel application itself when it writes the
# file). The user-visible, reproducible edits are: the active selection
# moved to L21, and every data column got a fraction wider (re-measured
# "best fit" widths), so that is what we replay here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (best-fit re-measurement) ------------------------------
# ColumnWidth is expressed in "characters" of the Normal style font and the
# engine (like Excel itself) snaps it to whole pixels, i.e. to steps of
# 1/6 of a character here. Requesting (target - 5/6) lands in the pixel
# bucket that reproduces the target stored width as closely as that
# quantization allows.
$colWidths = @{
    1  = 8.7109375
    2  = 9.42578125
    3  = 15.7109375
    4  = 15.85546875
    5  = 24
    6  = 29.7109375
    7  = 32.85546875
    8  = 20.5703125
    9  = 23.85546875
    10 = 18
    11 = 28
    12 = 23.7109375
}

foreach ($col in $colWidths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $colWidths[$col] - (5/6)
}

# --- Selection -------------------------------------------------------------
$ws.Range("L21").Select()
